# "Generate Report for Handback"
#
# Fills in the handback-report columns (Latest Target File / Latest
# Handback File / Latest Handback DateTime) for both language sheets and
# flips the Status from "Ready for handoff" to "Handed back: in sync with
# en-US" everywhere it appears (Overview + both language tabs share the
# same string).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    Shows up on Overview!E2:F3 and on column C ("Status") of both the
#    zh-cn and de-de tabs.
# ---------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E2").Value = $newStatus
$ov.Range("F2").Value = $newStatus
$ov.Range("E3").Value = $newStatus
$ov.Range("F3").Value = $newStatus

$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = $newStatus
$zh.Range("C3").Value = $newStatus

$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = $newStatus
$de.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# 2) zh-cn sheet: fill Latest Target File (I), Latest Handback File (J)
#    and Latest Handback DateTime (K) for both data rows, and hyperlink
#    the new "Latest Target File" cells the same way column A is linked.
# ---------------------------------------------------------------------
$zh.Range("I2").Value = "4cc237eb-894e-46d4-a753-e8edc05a87a2.md"
$zh.Range("J2").Value = "4cc237eb-894e-46d4-a753-e8edc05a87a2.21bee339a012b50bf842e42ad59fa4c786adcc14.zh-cn.xlf"
$zh.Range("K2").Value = "2016-08-17 09:00:48"

$zh.Range("I3").Value = "f0496682-930e-4021-af04-f1d1ae7cd233.md"
$zh.Range("J3").Value = "f0496682-930e-4021-af04-f1d1ae7cd233.3fced64c4860805d582eb0b6cba1b647fe724050.zh-cn.xlf"
$zh.Range("K3").Value = "2016-08-17 09:00:48"

# ---------------------------------------------------------------------
# 3) de-de sheet: same fields.
# ---------------------------------------------------------------------
$de.Range("I2").Value = "4cc237eb-894e-46d4-a753-e8edc05a87a2.md"
$de.Range("J2").Value = "4cc237eb-894e-46d4-a753-e8edc05a87a2.21bee339a012b50bf842e42ad59fa4c786adcc14.de-de.xlf"
$de.Range("K2").Value = "2016-08-17 09:00:56"

$de.Range("I3").Value = "f0496682-930e-4021-af04-f1d1ae7cd233.md"
$de.Range("J3").Value = "f0496682-930e-4021-af04-f1d1ae7cd233.3fced64c4860805d582eb0b6cba1b647fe724050.de-de.xlf"
$de.Range("K3").Value = "2016-08-17 09:00:56"

# ---------------------------------------------------------------------
# 4) Hyperlink the new "Latest Target File" cells (I2/I3) to the same
#    GitHub blob URLs column A already links to. Rebuild each sheet's
#    hyperlink collection in row order (A2, I2, A3, I3) so the new
#    relationship ids interleave the same way the original workbook did.
# ---------------------------------------------------------------------
$url4cc = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/be445d18fc0a4596efd3df3a9e4757636a5e8025/e2e/4cc237eb-894e-46d4-a753-e8edc05a87a2.md"
$urlf04 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/be445d18fc0a4596efd3df3a9e4757636a5e8025/e2e/f0496682-930e-4021-af04-f1d1ae7cd233.md"

foreach ($ws in @($zh, $de)) {
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $url4cc, [System.Type]::Missing, [System.Type]::Missing, "4cc237eb-894e-46d4-a753-e8edc05a87a2.md")
    $ws.Hyperlinks.Add($ws.Range("I2"), $url4cc, [System.Type]::Missing, [System.Type]::Missing, "4cc237eb-894e-46d4-a753-e8edc05a87a2.md")
    $ws.Hyperlinks.Add($ws.Range("A3"), $urlf04, [System.Type]::Missing, [System.Type]::Missing, "f0496682-930e-4021-af04-f1d1ae7cd233.md")
    $ws.Hyperlinks.Add($ws.Range("I3"), $urlf04, [System.Type]::Missing, [System.Type]::Missing, "f0496682-930e-4021-af04-f1d1ae7cd233.md")
}

# ---------------------------------------------------------------------
# 5) Widen columns that now hold the longer filenames/status text.
#    (Excel itself quantises ColumnWidth to its internal character grid,
#    so these land on the nearest representable width.)
# ---------------------------------------------------------------------
$ov.Columns.Item(5).ColumnWidth = 29.9777047293527
$ov.Columns.Item(6).ColumnWidth = 29.9777047293527

foreach ($ws in @($zh, $de)) {
    $ws.Columns.Item(3).ColumnWidth = 29.9777047293527
    $ws.Columns.Item(9).ColumnWidth = 40
    $ws.Columns.Item(10).ColumnWidth = 40
}

Write-Output "Generate Report for Handback: done"
